$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2066.1667
$ws.Range("I19").Value = 1923.75
$ws.Range("J19").Value = 2351
$ws.Range("K19").Value = 1923.75
$ws.Range("L19").Value = 2351
$ws.Range("M19").Value = -1748.75
$ws.Range("N19").Value = -2701

$ws.Range("H33").Value = 598
$ws.Range("I33").Value = 598
$ws.Range("K33").Value = 598
$ws.Range("M33").Value = -369

$ws.Range("H76").Value = 4415
$ws.Range("I76").Value = 4415
$ws.Range("K76").Value = 4415
$ws.Range("M76").Value = -4100

$ws.Range("H79").Value = 4415
$ws.Range("I79").Value = 4415
$ws.Range("K79").Value = 4415
$ws.Range("M79").Value = -3323

$ws.Range("H125").Value = 9999
$ws.Range("I125").Value = 9999
$ws.Range("K125").Value = 89991
$ws.Range("M125").Value = -87531

$ws.Range("H132").Value = 5961.5
$ws.Range("I132").Value = 2949.6667
$ws.Range("K132").Value = 8849.000100000001
$ws.Range("M132").Value = -6319.000100000001

$ws.Range("H137").Value = 2735.4167
$ws.Range("J137").Value = 3639.8333
$ws.Range("L137").Value = 10919.4999
$ws.Range("N137").Value = -16019.4999

$ws.Range("H138").Value = 7427.7627
$ws.Range("J138").Value = 8021.6226
$ws.Range("L138").Value = 24064.8678
$ws.Range("N138").Value = -34344.8678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1100
$ws.Range("I2").Value = 1100
$ws.Range("K2").Value = 1100
$ws.Range("M2").Value = -987

$ws.Range("H45").Value = 2100.9167
$ws.Range("I45").Value = 2100.9167
$ws.Range("K45").Value = 2100.9167
$ws.Range("M45").Value = -1723.9167

$ws.Range("H61").Value = 4145.9375
$ws.Range("I61").Value = 4713.3335
$ws.Range("K61").Value = 4713.3335
$ws.Range("M61").Value = -4501.3335

$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492

$ws.Range("H116").Value = 1100
$ws.Range("I116").Value = 1100
$ws.Range("K116").Value = 1100
$ws.Range("M116").Value = 1194

$ws.Range("H122").Value = 7160.5
$ws.Range("I122").Value = 7524.2144
$ws.Range("K122").Value = 22572.6432
$ws.Range("M122").Value = -20122.6432

$ws.Range("H132").Value = 3683.8572
$ws.Range("I132").Value = 2697.5
$ws.Range("K132").Value = 8092.5
$ws.Range("M132").Value = -5562.5

$ws.Range("H136").Value = 4145.9375
$ws.Range("I136").Value = 4713.3335
$ws.Range("K136").Value = 14140.0005
$ws.Range("M136").Value = -11590.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1100
$ws.Range("I3").Value = 1100
$ws.Range("K3").Value = 1100
$ws.Range("M3").Value = -986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1596.3334
$ws.Range("I31").Value = 1303.5385
$ws.Range("J31").Value = 3499.5
$ws.Range("K31").Value = 1303.5385
$ws.Range("L31").Value = 3499.5
$ws.Range("M31").Value = -1008.5385
$ws.Range("N31").Value = -4089.5

$ws.Range("H34").Value = 1596.3334
$ws.Range("I34").Value = 1303.5385
$ws.Range("J34").Value = 3499.5
$ws.Range("K34").Value = 1303.5385
$ws.Range("L34").Value = 3499.5
$ws.Range("M34").Value = -1101.5385
$ws.Range("N34").Value = -3903.5

$ws.Range("H86").Value = 50278.5
$ws.Range("I86").Value = 11741.333
$ws.Range("K86").Value = 11741.333
$ws.Range("M86").Value = -10618.333

$ws.Range("H89").Value = 50278.5
$ws.Range("I89").Value = 11741.333
$ws.Range("K89").Value = 58706.665
$ws.Range("M89").Value = -53090.665

$ws.Range("H96").Value = 55333.332
$ws.Range("J96").Value = 55333.332
$ws.Range("L96").Value = 55333.332
$ws.Range("N96").Value = -60825.332

$ws.Range("H132").Value = 3832.2942
$ws.Range("I132").Value = 2407.6667
$ws.Range("K132").Value = 7223.000100000001
$ws.Range("M132").Value = -4693.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 87.5
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -976

$ws.Range("H4").Value = 100597130
$ws.Range("I4").Value = 17417470
$ws.Range("J4").Value = 266956450
$ws.Range("K4").Value = 52252410
$ws.Range("L4").Value = 800869350
$ws.Range("M4").Value = -52252298
$ws.Range("N4").Value = -800869574

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H18").Value = 400
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H122").Value = 854.5
$ws.Range("J122").Value = 1022
$ws.Range("L122").Value = 9198
$ws.Range("N122").Value = -14098

$ws.Range("H131").Value = 1662.25
$ws.Range("I131").Value = 1624.5
$ws.Range("K131").Value = 4873.5
$ws.Range("M131").Value = 166.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1698.0714
$ws.Range("I97").Value = 314.5
$ws.Range("K97").Value = 314.5
$ws.Range("M97").Value = 181.5

$ws.Range("H102").Value = 4266.6665
$ws.Range("I102").Value = 4266.6665
$ws.Range("K102").Value = 4266.6665
$ws.Range("M102").Value = -2644.6665

$ws.Range("H132").Value = 9273.75
$ws.Range("I132").Value = 12548.5
$ws.Range("K132").Value = 37645.5
$ws.Range("M132").Value = -35115.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2722.2144
$ws.Range("I16").Value = 2413.182
$ws.Range("K16").Value = 2413.182
$ws.Range("M16").Value = -2243.182

$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214

$ws.Range("H40").Value = 3628.4285
$ws.Range("I40").Value = 3616.5
$ws.Range("J40").Value = 3700
$ws.Range("K40").Value = 3616.5
$ws.Range("L40").Value = 3700
$ws.Range("M40").Value = -3480.5
$ws.Range("N40").Value = -3972

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws.Range("H68").Value = 2998.9167
$ws.Range("I68").Value = 2998.9167
$ws.Range("K68").Value = 2998.9167
$ws.Range("M68").Value = -2249.9167

$ws.Range("H71").Value = 2998.9167
$ws.Range("I71").Value = 2998.9167
$ws.Range("K71").Value = 14994.5835
$ws.Range("M71").Value = -11250.5835

$ws.Range("H132").Value = 3281.7727
$ws.Range("I132").Value = 2679.077
$ws.Range("J132").Value = 4152.3335
$ws.Range("K132").Value = 8037.231000000001
$ws.Range("L132").Value = 12457.0005
$ws.Range("M132").Value = -5507.231000000001
$ws.Range("N132").Value = -17517.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3128
$ws.Range("I62").Value = 3242
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 3242
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -2618
$ws.Range("N62").Value = -4148

$ws.Range("H65").Value = 3128
$ws.Range("I65").Value = 3242
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 16210
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -13090
$ws.Range("N65").Value = -20740

$ws.Range("H100").Value = 995.3333
$ws.Range("I100").Value = 1247.5
$ws.Range("J100").Value = 491
$ws.Range("K100").Value = 2495
$ws.Range("L100").Value = 982
$ws.Range("M100").Value = -1954
$ws.Range("N100").Value = -2064

$ws.Range("H122").Value = 3948.1667
$ws.Range("I122").Value = 3923.5
$ws.Range("K122").Value = 11770.5
$ws.Range("M122").Value = -9320.5

$ws.Range("H126").Value = 2459.8
$ws.Range("I126").Value = 2199.75
$ws.Range("K126").Value = 6599.25
$ws.Range("M126").Value = -4129.25

$ws.Range("H132").Value = 2803.3462
$ws.Range("I132").Value = 2324.0588
$ws.Range("K132").Value = 6972.176399999999
$ws.Range("M132").Value = -4442.176399999999

